$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address line into two paragraphs:
#    "19964 Lorena Circle, Castro Valley CA 94546"
#    becomes
#    "19964 Lorena Circle"
#    "Castro Valley, CA 94546"
$addrRng = $d.Content.Duplicate
$addrRng.Find.Execute("19964 Lorena Circle, Castro Valley CA 94546")
$addrRng.Text = "19964 Lorena Circle"
$addrRng.Collapse(0)          # wdCollapseEnd
$addrRng.InsertParagraphAfter()
$addrRng.Collapse(0)          # wdCollapseEnd
$addrRng.MoveStart(1, 1)      # wdCharacter - step past the new paragraph mark
$addrRng.InsertAfter("Castro Valley, CA 94546")

# 3. Remove the blank paragraph that immediately follows "Board of Directors".
$bodRng = $d.Content.Duplicate
$bodRng.Find.Execute("Board of Directors")
$bodPara = $bodRng.Paragraphs(1)
$blankPara = $bodPara.Next()
$blankPara.Range.Delete()
